# Add two new expense rows below the existing data (rows 2-3), extending
# the sheet from A1:E3 to A1:E5, matching the "Food"/"Transport" entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Food"
$ws.Range("C4").Value = 500
$ws.Range("D4").Value = "July 10, 2026"
# A lone leading apostrophe is Excel's "treat as text" marker; it yields an
# empty text cell (like the existing blank Description cells E2/E3) rather
# than clearing the cell outright (which is what Value = "" would do).
$ws.Range("E4").Value = "'"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Transport"
$ws.Range("C5").Value = 1111
$ws.Range("D5").Value = "May 18 2025"
$ws.Range("E5").Value = "'"
